$wb = $excel.ActiveWorkbook

# Rename the sole data sheet from "DOE21E" to the generic "Sheet1".
$ws = $wb.Worksheets.Item("DOE21E")
$ws.Name = "Sheet1"

# Remove the (unused) external workbook reference entirely - converts
# the external link formulas/definitions to static values and drops the
# <externalReferences> element + externalLink1.xml part on save.
$sources = $wb.LinkSources()
if ($sources) {
    foreach ($src in $sources) {
        $wb.BreakLink($src, 1)
    }
}
